$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.745.47"
$ws.Range("E2").Value = "  -2.22%  "
$ws.Range("D3").Value = "2.339.87"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'502.50"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").Value = "'128.51"
$ws.Range("E6").Value = "  -3.51%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.535"
$ws.Range("E8").Value = "  -3.36%  "
$ws.Range("D9").Value = "2.346.95"
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("D10").Value = "'0.0976"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "'4.78"
$ws.Range("E12").Value = "  +3.00%  "
$ws.Range("D13").Value = "'0.320"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").Value = "2.753.85"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("D15").Value = "55.715.49"
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").Value = "'21.61"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "2.339.38"
$ws.Range("E18").Value = "  -3.01%  "
$ws.Range("D19").Value = "'9.94"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("D20").Value = "'309.31"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").Value = "'6.19"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  -4.05%  "
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("E27").Value = "  -3.88%  "
$ws.Range("D28").Value = "'7.10"
$ws.Range("E28").Value = "  -4.62%  "
$ws.Range("D29").Value = "'172.32"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("E30").Value = "  -2.73%  "
$ws.Range("D31").Value = "0.0₃0700"
$ws.Range("E31").Value = "  -3.80%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'5.77"
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("E35").Value = "  -5.74%  "
$ws.Range("D36").Value = "'17.61"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'3.65"
$ws.Range("E38").Value = "  -5.04%  "
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").Value = "'0.821"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "'36.09"
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("E41").Value = "  -4.42%  "
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").Value = "'126.74"
$ws.Range("E43").Value = "  -4.47%  "
$ws.Range("D44").Value = "'4.71"
$ws.Range("E44").Value = "  -5.53%  "
$ws.Range("D45").Value = "'0.554"
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").Value = "'236.47"
$ws.Range("E47").Value = "  -6.08%  "
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").Value = "'0.954"
$ws.Range("E51").Value = "  +0.25%  "
